$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column F entirely (it held "RequestProcessingType"), shifting
# everything from column G onward one column to the left.
$ws.Columns("F:F").Delete()

# Rename headers to their standardized field names (row 1, after the shift).
$ws.Range("F1").Value = "ResponseByFields"
$ws.Range("I1").Value = "Action"
$ws.Range("J1").Value = "ExcludeFields"
$ws.Range("K1").Value = "StatusCode"
$ws.Range("N1").Value = "Security"
$ws.Range("O1").Value = "Tags"

# Restore the cursor/selection to the cell where the last column used to be.
$ws.Range("P1").Select() | Out-Null
